$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Tipo" column (D) so the
# existing D1/D2 (Tipo/single) shift to E1/E2, and the new MAE values
# occupy the now-empty D column.
$ws.Columns.Item(4).Insert()

# Header for the new column (match the existing header formatting:
# bold font, thin border, centered horizontal/top vertical alignment)
$ws.Range("D1").Value = "MAE"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Borders.LineStyle = 1

# New MAE value for row 2
$ws.Range("D2").Value = 0.5256692707344185

# Minor correction to the existing MSE value in B2
$ws.Range("B2").Value = 0.4108830721519877
